$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 4; this shifts existing rows 4..96 down to 5..97
$ws.Rows.Item(4).Insert()

# Fill in the new row 4 with a new weekly price observation.
# Columns A, B, C, E, F, G, H, I, L, N, O, Q, R keep the same values as the
# (now shifted) row below it; only D (date), J, K, M and P change.
$ws.Cells.Item(4, 1).Value = 1
$ws.Cells.Item(4, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(4, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(4, 4).Value = "9/21/2023"
$ws.Cells.Item(4, 5).Value = 15
$ws.Cells.Item(4, 6).Value = 100112031
$ws.Cells.Item(4, 7).Value = "Poroto verde"
$ws.Cells.Item(4, 8).Value = "Sin especificar"
$ws.Cells.Item(4, 9).Value = "Primera"
$ws.Cells.Item(4, 10).Value = 700
$ws.Cells.Item(4, 11).Value = 900
$ws.Cells.Item(4, 12).Value = 1000
$ws.Cells.Item(4, 13).Value = 943
$ws.Cells.Item(4, 14).Value = "`$/kilo"
$ws.Cells.Item(4, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(4, 16).Value = 943
$ws.Cells.Item(4, 17).Value = 1
$ws.Cells.Item(4, 18).Value = "Hortaliza"
